$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 15.79399266666667
$ws.Cells.Item(2, 8).Value = 47.381978
$ws.Cells.Item(2, 9).Value = 0.2968109173698557
$ws.Cells.Item(2, 10).Value = 0.2968109173698557
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 2.618716333333334
$ws.Cells.Item(2, 14).Value = 7.856149000000001
$ws.Cells.Item(2, 15).Value = 0.07115908183301342
$ws.Cells.Item(2, 16).Value = 0.07115908183301341
$ws.Cells.Item(2, 17).Value = 41.3599865647469
$ws.Cells.Item(2, 18).Value = 372.2398790827221
$ws.Cells.Item(2, 19).Value = 0.02112079235805335
$ws.Cells.Item(2, 20).Value = 0.02112079235805334

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 15.79399266666667
$ws.Cells.Item(3, 8).Value = 47.381978
$ws.Cells.Item(3, 9).Value = 0.2968109173698557
$ws.Cells.Item(3, 10).Value = 0.2968109173698557
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 15.503283
$ws.Cells.Item(3, 14).Value = 46.509849
$ws.Cells.Item(3, 15).Value = 0.4212748702999519
$ws.Cells.Item(3, 16).Value = 0.4212748702999519
$ws.Cells.Item(3, 17).Value = 244.8587380112581
$ws.Cells.Item(3, 18).Value = 2203.728642101322
$ws.Cells.Item(3, 19).Value = 0.1250389807185957
$ws.Cells.Item(3, 20).Value = 0.1250389807185957

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 15.79399266666667
$ws.Cells.Item(4, 8).Value = 47.381978
$ws.Cells.Item(4, 9).Value = 0.2968109173698557
$ws.Cells.Item(4, 10).Value = 0.2968109173698557
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 18.67887366666666
$ws.Cells.Item(4, 14).Value = 56.036621
$ws.Cells.Item(4, 15).Value = 0.5075660478670347
$ws.Cells.Item(4, 16).Value = 0.5075660478670347
$ws.Cells.Item(4, 17).Value = 295.0139937129264
$ws.Cells.Item(4, 18).Value = 2655.125943416338
$ws.Cells.Item(4, 19).Value = 0.1506511442932067
$ws.Cells.Item(4, 20).Value = 0.1506511442932066

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 24.86954866666666
$ws.Cells.Item(5, 8).Value = 74.60864599999999
$ws.Cells.Item(5, 9).Value = 0.4673646309781075
$ws.Cells.Item(5, 10).Value = 0.4673646309781075
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 2.618716333333334
$ws.Cells.Item(5, 14).Value = 7.856149000000001
$ws.Cells.Item(5, 15).Value = 0.07115908183301342
$ws.Cells.Item(5, 16).Value = 0.07115908183301341
$ws.Cells.Item(5, 17).Value = 65.12629329602822
$ws.Cells.Item(5, 18).Value = 586.1366396642541
$ws.Cells.Item(5, 19).Value = 0.03325723802162727
$ws.Cells.Item(5, 20).Value = 0.03325723802162727

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 24.86954866666666
$ws.Cells.Item(6, 8).Value = 74.60864599999999
$ws.Cells.Item(6, 9).Value = 0.4673646309781075
$ws.Cells.Item(6, 10).Value = 0.4673646309781075
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 15.503283
$ws.Cells.Item(6, 14).Value = 46.509849
$ws.Cells.Item(6, 15).Value = 0.4212748702999519
$ws.Cells.Item(6, 16).Value = 0.4212748702999519
$ws.Cells.Item(6, 17).Value = 385.559651061606
$ws.Cells.Item(6, 18).Value = 3470.036859554454
$ws.Cells.Item(6, 19).Value = 0.1968889742980871
$ws.Cells.Item(6, 20).Value = 0.1968889742980871

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 24.86954866666666
$ws.Cells.Item(7, 8).Value = 74.60864599999999
$ws.Cells.Item(7, 9).Value = 0.4673646309781075
$ws.Cells.Item(7, 10).Value = 0.4673646309781075
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 18.67887366666666
$ws.Cells.Item(7, 14).Value = 56.036621
$ws.Cells.Item(7, 15).Value = 0.5075660478670347
$ws.Cells.Item(7, 16).Value = 0.5075660478670347
$ws.Cells.Item(7, 17).Value = 464.535157691685
$ws.Cells.Item(7, 18).Value = 4180.816419225165
$ws.Cells.Item(7, 19).Value = 0.2372184186583932
$ws.Cells.Item(7, 20).Value = 0.2372184186583932

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 12.54876233333333
$ws.Cells.Item(8, 8).Value = 37.646287
$ws.Cells.Item(8, 9).Value = 0.2358244516520368
$ws.Cells.Item(8, 10).Value = 0.2358244516520368
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 2.618716333333334
$ws.Cells.Item(8, 14).Value = 7.856149000000001
$ws.Cells.Item(8, 15).Value = 0.07115908183301342
$ws.Cells.Item(8, 16).Value = 0.07115908183301341
$ws.Cells.Item(8, 17).Value = 32.86164888541812
$ws.Cells.Item(8, 18).Value = 295.7548399687631
$ws.Cells.Item(8, 19).Value = 0.01678105145333281
$ws.Cells.Item(8, 20).Value = 0.0167810514533328

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 12.54876233333333
$ws.Cells.Item(9, 8).Value = 37.646287
$ws.Cells.Item(9, 9).Value = 0.2358244516520368
$ws.Cells.Item(9, 10).Value = 0.2358244516520368
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 15.503283
$ws.Cells.Item(9, 14).Value = 46.509849
$ws.Cells.Item(9, 15).Value = 0.4212748702999519
$ws.Cells.Item(9, 16).Value = 0.4212748702999519
$ws.Cells.Item(9, 17).Value = 194.547013753407
$ws.Cells.Item(9, 18).Value = 1750.923123780663
$ws.Cells.Item(9, 19).Value = 0.0993469152832691
$ws.Cells.Item(9, 20).Value = 0.09934691528326907

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 12.54876233333333
$ws.Cells.Item(10, 8).Value = 37.646287
$ws.Cells.Item(10, 9).Value = 0.2358244516520368
$ws.Cells.Item(10, 10).Value = 0.2358244516520368
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 18.67887366666666
$ws.Cells.Item(10, 14).Value = 56.036621
$ws.Cells.Item(10, 15).Value = 0.5075660478670347
$ws.Cells.Item(10, 16).Value = 0.5075660478670347
$ws.Cells.Item(10, 17).Value = 234.3967462973585
$ws.Cells.Item(10, 18).Value = 2109.570716676227
$ws.Cells.Item(10, 19).Value = 0.1196964849154349
$ws.Cells.Item(10, 20).Value = 0.1196964849154349
